# Prepare for restud submission
# Restructure the header of the hit_miss_rule table: split the single header
# row (row 4) into two rows - a new top row 3 with the main column titles
# ("Rule" / "% incorrectly " / "% incorrectly " / "Overall Error Rate") and a
# second line in row 4 with the continuation text ("assigned to control " /
# " assigned to treatment"). The H:I columns (TI / T2) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Build new row 3 with the formatting currently used by row 4 (C4:F4) ---
$ws.Range("C4:F4").Copy()
$ws.Range("C3:F3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = 0

$ws.Range("C3").Value2 = "Rule"
$ws.Range("D3").Value2 = "% incorrectly "
$ws.Range("E3").Value2 = "% incorrectly "
$ws.Range("F3").Value2 = "Overall Error Rate"

# remove the bottom (double) border on the new top row - it now only keeps
# the thin top border
$ws.Range("C3:F3").Borders(9).LineStyle = -4142   # xlEdgeBottom / xlLineStyleNone

# --- 2. Turn old row 4 into the second header line, re-using the bottom
#        double-border / centred style already used by row 10 (C10:F10) ---
$ws.Range("D10").Copy()
$ws.Range("C4:F4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = 0

$ws.Range("C4").ClearContents()
$ws.Range("D4").Value2 = "assigned to control "
$ws.Range("E4").Value2 = " assigned to treatment"
$ws.Range("F4").ClearContents()

# --- 3. Row 10 label: "Allow choice" text is unchanged, only its shared
#        string slot shifts internally - no explicit action required beyond
#        making sure the text itself still reads "Allow choice" ---
$ws.Range("C10").Value2 = "Allow choice"

# --- 4. Column widths - the text is now shorter (wrapped onto two lines) so
#        the best-fit column widths shrink (values chosen so the engine's
#        pixel-quantised ColumnWidth lands as close as possible to the
#        target best-fit widths of 17.1796875 / 19.6328125 / 15.7265625) ---
$ws.Columns("D").ColumnWidth = 16.333333333333332
$ws.Columns("E").ColumnWidth = 18.833333333333332
$ws.Columns("F").ColumnWidth = 14.833333333333334

# --- 5. Selection / active cell now starts at the new top-left header cell ---
[void]$ws.Range("C3:F10").Select()
